$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 323, pushing all following rows
# (old 323-377) down to 326-380.
$ws.Rows("323:325").Insert()

# New row 323: Especial / 108 / 18000 / $/caja 10 unidades / 1800 / 10
$ws.Range("A323").Value = 3
$ws.Range("B323").Value = "Femacal de La Calera"
$ws.Range("C323").Value = "Coquimbo"
$ws.Range("D323").Value = 44505
$ws.Range("E323").Value = 5
$ws.Range("F323").Value = "Fruta"
$ws.Range("G323").Value = 100108
$ws.Range("H323").Value = "Tropicales y subtropicales"
$ws.Range("I323").Value = 100108005
$ws.Range("J323").Value = "Piña"
$ws.Range("K323").Value = "Caramelo"
$ws.Range("L323").Value = "Especial"
$ws.Range("M323").Value = 108
$ws.Range("N323").Value = 18000
$ws.Range("O323").Value = 18000
$ws.Range("P323").Value = 18000
$ws.Range("Q323").Value = "$/caja 10 unidades"
$ws.Range("R323").Value = "Ecuador"
$ws.Range("S323").Value = 1800
$ws.Range("T323").Value = 10

# New row 324: Primera / 162 / 18000 / $/caja 12 unidades / 1500 / 12
$ws.Range("A324").Value = 3
$ws.Range("B324").Value = "Femacal de La Calera"
$ws.Range("C324").Value = "Coquimbo"
$ws.Range("D324").Value = 44505
$ws.Range("E324").Value = 5
$ws.Range("F324").Value = "Fruta"
$ws.Range("G324").Value = 100108
$ws.Range("H324").Value = "Tropicales y subtropicales"
$ws.Range("I324").Value = 100108005
$ws.Range("J324").Value = "Piña"
$ws.Range("K324").Value = "Caramelo"
$ws.Range("L324").Value = "Primera"
$ws.Range("M324").Value = 162
$ws.Range("N324").Value = 18000
$ws.Range("O324").Value = 18000
$ws.Range("P324").Value = 18000
$ws.Range("Q324").Value = "$/caja 12 unidades"
$ws.Range("R324").Value = "Ecuador"
$ws.Range("S324").Value = 1500
$ws.Range("T324").Value = 12

# New row 325: Segunda / 162 / 18000 / $/caja 14 unidades / 1286 / 14
$ws.Range("A325").Value = 3
$ws.Range("B325").Value = "Femacal de La Calera"
$ws.Range("C325").Value = "Coquimbo"
$ws.Range("D325").Value = 44505
$ws.Range("E325").Value = 5
$ws.Range("F325").Value = "Fruta"
$ws.Range("G325").Value = 100108
$ws.Range("H325").Value = "Tropicales y subtropicales"
$ws.Range("I325").Value = 100108005
$ws.Range("J325").Value = "Piña"
$ws.Range("K325").Value = "Caramelo"
$ws.Range("L325").Value = "Segunda"
$ws.Range("M325").Value = 162
$ws.Range("N325").Value = 18000
$ws.Range("O325").Value = 18000
$ws.Range("P325").Value = 18000
$ws.Range("Q325").Value = "$/caja 14 unidades"
$ws.Range("R325").Value = "Ecuador"
$ws.Range("S325").Value = 1286
$ws.Range("T325").Value = 14
